$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: clear the effort value for H3 (keep style)
$ws.Range("H3").ClearContents()

# Row 4: update percentage and effort columns
$ws.Range("G4").Value = 33.333
$ws.Range("H4").Value = 5

# Row 5: update percentage and effort columns
$ws.Range("G5").Value = 33.333
$ws.Range("H5").Value = 5

# Row 6: clear the "X" marks and percentage, set effort to 0
$ws.Range("B6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 0

# Row 7: update percentage and effort columns
$ws.Range("G7").Value = 33.333
$ws.Range("H7").Value = 5

# Row 8: set total percentage
$ws.Range("G8").Value = 100

# Update the active selection
$ws.Range("L16").Select() | Out-Null
